# fix(publipostage): Correct status name
#
# Renames the "statut_label" value "bleu" -> "noir", and rewords the four
# "statut_name" descriptions to replace "... et / ou publication posté ..."
# with "... postés ou publiés ..." phrasing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "bleu"                                                  = "noir"
    "pas de résultat ni de publication"                     = "pas de résultat postés ni publiés"
    "résultat et / ou publication posté"                    = "résultat postés ou publiés"
    "résultat et / ou publication posté dans les 36 mois"    = "résultat postés ou publiés dans les 36 mois"
    "résultat et / ou publication posté dans les 12 mois"    = "résultat postés ou publiés dans les 12 mois"
}

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -is [string] -and $replacements.ContainsKey($val)) {
            $cell.Value = $replacements[$val]
        }
    }
}
